$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.13210666179657
$ws.Range("B1").Value = 2.218048810958862
$ws.Range("C1").Value = 10.76741027832031
$ws.Range("D1").Value = 2.275628328323364
$ws.Range("E1").Value = 1.283723115921021
